$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# Insert a new column before F (shifts old F/G/H -> G/H/I)
$ws.Columns("F:F").Insert()

Write-Output "done"
